$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.101.71'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.787.54'
$ws.Range("E3").Value = '  -2.82%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.27'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  -1.69%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.86'
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0709'
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.045.98'
$ws.Range("E12").Value = '  -2.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.790.62'
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.86'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  -3.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.048.44'
$ws.Range("E16").Value = '  -1.20%  '
$ws.Range("E17").Value = '  -4.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.86'
$ws.Range("E18").Value = '  -2.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.43'
$ws.Range("E19").Value = '  -2.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0790'
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.81'
$ws.Range("E22").Value = '  -4.94%  '
$ws.Range("E23").Value = '  -4.54%  '
$ws.Range("E24").Value = '  -2.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.38'
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.35'
$ws.Range("E27").Value = '  -2.83%  '
$ws.Range("E28").Value = '  -3.04%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -4.40%  '
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("E34").Value = '  -5.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.396.46'
$ws.Range("E35").Value = '  -4.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.644'
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0187'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.21'
$ws.Range("E39").Value = '  +2.77%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  -5.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '78.06'
$ws.Range("E43").Value = '  -5.24%  '
$ws.Range("E44").Value = '  +16.33%  '
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.53'
$ws.Range("E46").Value = '  +3.16%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0499'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '108.14'
$ws.Range("E48").Value = '  +1.35%  '
$ws.Range("E49").Value = '  -3.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.945.40'
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.34%  '
